$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 28.666666
$ws.Range("I8").Value = 28.666666
$ws.Range("K8").Value = 85.99999800000001
$ws.Range("M8").Value = 53.00000199999999

$ws.Range("H9").Value = 143.28572
$ws.Range("I9").Value = 163.5
$ws.Range("J9").Value = 92.75
$ws.Range("K9").Value = 163.5
$ws.Range("L9").Value = 92.75
$ws.Range("M9").Value = 5.5
$ws.Range("N9").Value = -430.75

$ws.Range("H64").Value = 4399.421
$ws.Range("I64").Value = 5148.2856
$ws.Range("K64").Value = 5148.2856
$ws.Range("M64").Value = -4900.2856

$ws.Range("H67").Value = 4399.421
$ws.Range("I67").Value = 5148.2856
$ws.Range("K67").Value = 5148.2856
$ws.Range("M67").Value = -4290.2856

$ws.Range("H86").Value = 22857.572
$ws.Range("I86").Value = 20001
$ws.Range("K86").Value = 20001
$ws.Range("M86").Value = -18878

$ws.Range("H89").Value = 22857.572
$ws.Range("I89").Value = 20001
$ws.Range("K89").Value = 100005
$ws.Range("M89").Value = -94389

$ws.Range("H113").Value = 3250.3333
$ws.Range("I113").Value = 3219.125
$ws.Range("K113").Value = 3219.125
$ws.Range("M113").Value = 34.875

$ws.Range("H131").Value = 2157
$ws.Range("I131").Value = 2539.8
$ws.Range("K131").Value = 7619.400000000001
$ws.Range("M131").Value = -2579.400000000001

$ws.Range("H132").Value = 17959.334
$ws.Range("I132").Value = 17959.334
$ws.Range("K132").Value = 53878.00199999999
$ws.Range("M132").Value = -51348.00199999999

$ws.Range("H138").Value = 3717.5366
$ws.Range("I138").Value = 4105.9414
$ws.Range("J138").Value = 3442.4167
$ws.Range("K138").Value = 12317.8242
$ws.Range("L138").Value = 10327.2501
$ws.Range("M138").Value = -7177.824199999999
$ws.Range("N138").Value = -20607.2501

$ws.Range("H141").Value = 5769.926
$ws.Range("I141").Value = 3119.65
$ws.Range("K141").Value = 9358.950000000001
$ws.Range("M141").Value = -4178.950000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2999.1428

$ws.Range("H32").Value = 183317.6
$ws.Range("I32").Value = 184832.45
$ws.Range("K32").Value = 184832.45
$ws.Range("M32").Value = -184545.45

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 4833.385
$ws.Range("I132").Value = 3104.1667
$ws.Range("J132").Value = 7991.087
$ws.Range("K132").Value = 9312.500100000001
$ws.Range("L132").Value = 23973.261
$ws.Range("M132").Value = -6782.500100000001
$ws.Range("N132").Value = -29033.261

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 283.33334
$ws.Range("I22").Value = 283.33334
$ws.Range("K22").Value = 283.33334
$ws.Range("M22").Value = -110.33334

$ws.Range("H99").Value = 10780.909
$ws.Range("I99").Value = 15914.143
$ws.Range("J99").Value = 1797.75
$ws.Range("K99").Value = 15914.143
$ws.Range("L99").Value = 1797.75
$ws.Range("M99").Value = -14416.143
$ws.Range("N99").Value = -4793.75

$ws.Range("H132").Value = 199998
$ws.Range("J132").Value = 199998
$ws.Range("L132").Value = 199998
$ws.Range("N132").Value = -210118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8893.6
$ws.Range("I16").Value = 8893.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 8893.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -8606.6
$ws.Range("N16").ClearContents()

$ws.Range("H53").Value = 44187
$ws.Range("J53").Value = 40670.75
$ws.Range("L53").Value = 40670.75
$ws.Range("N53").Value = -41884.75

$ws.Range("H58").Value = 4738.12
$ws.Range("I58").Value = 3158.95
$ws.Range("K58").Value = 3158.95
$ws.Range("M58").Value = -2955.95

$ws.Range("H113").Value = 8893.6
$ws.Range("I113").Value = 8893.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8893.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6723.6
$ws.Range("N113").ClearContents()

$ws.Range("H115").Value = 54945
$ws.Range("J115").Value = 54945
$ws.Range("L115").Value = 54945
$ws.Range("N115").Value = -57295

$ws.Range("H136").Value = 4738.12
$ws.Range("I136").Value = 3158.95
$ws.Range("K136").Value = 9476.849999999999
$ws.Range("M136").Value = -6926.849999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 912
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 912
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5472
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -5698

$ws.Range("H12").Value = 359.14285
$ws.Range("I12").Value = 575
$ws.Range("J12").Value = 323.16666
$ws.Range("K12").Value = 1725
$ws.Range("L12").Value = 969.4999799999999
$ws.Range("M12").Value = -1552
$ws.Range("N12").Value = -1315.49998

$ws.Range("H122").Value = 1468055
$ws.Range("J122").Value = 2160.889
$ws.Range("L122").Value = 19448.001
$ws.Range("N122").Value = -24348.001

$ws.Range("H131").Value = 2094.6184
$ws.Range("I131").Value = 809.2222
$ws.Range("J131").Value = 2267.2837
$ws.Range("K131").Value = 2427.6666
$ws.Range("L131").Value = 6801.8511
$ws.Range("M131").Value = 2612.3334
$ws.Range("N131").Value = -16881.8511

$ws.Range("H132").Value = 2849.4285
$ws.Range("I132").Value = 3024.3333
$ws.Range("K132").Value = 27218.9997
$ws.Range("M132").Value = -24688.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2729.353
$ws.Range("I126").Value = 2391.25
$ws.Range("J126").Value = 3029.889
$ws.Range("K126").Value = 7173.75
$ws.Range("L126").Value = 9089.667000000001
$ws.Range("M126").Value = -4703.75
$ws.Range("N126").Value = -14029.667

$ws.Range("H141").Value = 73166.71000000001
$ws.Range("J141").Value = 73166.71000000001
$ws.Range("L141").Value = 73166.71000000001
$ws.Range("N141").Value = -83526.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3998.1936
$ws.Range("I46").Value = 1852.5714
$ws.Range("J46").Value = 4624
$ws.Range("K46").Value = 1852.5714
$ws.Range("L46").Value = 4624
$ws.Range("M46").Value = -1664.5714
$ws.Range("N46").Value = -5000

$ws.Range("H93").Value = 4106.3335
$ws.Range("I93").Value = 1250
$ws.Range("J93").Value = 5534.5
$ws.Range("K93").Value = 1250
$ws.Range("L93").Value = 5534.5
$ws.Range("M93").Value = -2
$ws.Range("N93").Value = -8030.5

$ws.Range("H94").Value = 89999.5
$ws.Range("J94").Value = 89999.5
$ws.Range("L94").Value = 89999.5
$ws.Range("N94").Value = -91351.5

$ws.Range("H132").Value = 3346.4285
$ws.Range("J132").Value = 3897.6667
$ws.Range("L132").Value = 11693.0001
$ws.Range("N132").Value = -16753.0001

$ws.Range("H133").Value = 84156.5
$ws.Range("J133").Value = 88776.664
$ws.Range("L133").Value = 88776.664
$ws.Range("N133").Value = -93836.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 69976.336
$ws.Range("J127").Value = 69976.336
$ws.Range("L127").Value = 69976.336
$ws.Range("N127").Value = -79896.336

$ws.Range("H128").Value = 52856.285
$ws.Range("J128").Value = 56665.832
$ws.Range("L128").Value = 56665.832
$ws.Range("N128").Value = -66625.83199999999

$ws.Range("H132").Value = 2523.45
$ws.Range("I132").Value = 1897.7858
$ws.Range("K132").Value = 5693.357400000001
$ws.Range("M132").Value = -3163.357400000001
